# Generate Report for Handoff
# For the four in-progress localization files (0708a165, 28be8992, 541327b8,
# fd604e61) a new handoff xliff was generated: the Priority is promoted from
# "low" to "ht", and the "Latest Handoff Datetime" is bumped to the new
# generation timestamp, on both the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zh = $wb.Worksheets.Item("zh-cn")
$de = $wb.Worksheets.Item("de-de")

# zh-cn sheet: rows 4-7 -> Priority (E) "low" -> "ht"; Handoff Datetime (H) refreshed
foreach ($r in 4..7) {
    $zh.Cells.Item($r, 5).Value = "ht"
    $zh.Cells.Item($r, 8).Value = "2016-08-30 06:33:28"
}

# de-de sheet: rows 4-7 -> Priority (E) "low" -> "ht"; Handoff Datetime (H) refreshed
foreach ($r in 4..7) {
    $de.Cells.Item($r, 5).Value = "ht"
    $de.Cells.Item($r, 8).Value = "2016-08-30 06:33:33"
}

# Overview sheet: "Latest HO Xliff Generate Date" (G) mirrors the de-de
# handoff datetime for these same four files, so it moves in lockstep.
foreach ($r in 4..7) {
    $overview.Cells.Item($r, 7).Value = "2016-08-30 06:33:33"
}
